# Applies the "Jupyter notebook used as an interface for users to use
# co2mpas driver library" edit to sample.xlsx:
#   - inputs!B2 (desired_velocity):   10    -> 40
#   - config!B2 (vehicle_id):         35135 -> 39393
#   - time_series: new column B (velocities), B2=5 then each subsequent
#     cell = previous + 5, using a literal then a normal formula then a
#     shared formula spanning B4:B22 (mirrors column A's layout)
#   - restore each sheet's remembered selection / active cell, and make
#     "inputs" the active (tabSelected) sheet.

$wb = $excel.ActiveWorkbook

$wsInputs = $wb.Worksheets.Item("inputs")
$wsConfig = $wb.Worksheets.Item("config")
$wsVehicle = $wb.Worksheets.Item("vehicle_inputs")
$wsTime = $wb.Worksheets.Item("time_series")

# --- data edits -----------------------------------------------------

$wsInputs.Range("B2").Value = 40
$wsConfig.Range("B2").Value = 39393

$wsTime.Range("B2").Value = 5
$wsTime.Range("B3").Formula = "=B2 + 5"
$wsTime.Range("B4:B22").Formula = "=B3 + 5"

# --- view / selection state ------------------------------------------
# Select in the same order the sheets appear so that the sheet selected
# last ("inputs") ends up the active / tabSelected sheet, matching the
# target file.

$wsConfig.Range("C11").Select()
$wsVehicle.Range("F20").Select()
$wsTime.Range("F15").Select()
$wsInputs.Range("A4").Select()
